$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(6, "Bernadette", "Leclerc", "bernadette.leclerc@example.com", "female", "Switzerland", "2025-03-01 15:01:58"),
    @(7, "Zoe", "Lavigne", "zoe.lavigne@example.com", "female", "Canada", "2025-03-01 15:01:58"),
    @(8, "Nash", "Westerik", "nash.westerik@example.com", "male", "Netherlands", "2025-03-01 15:01:58"),
    @(9, "Abbas", "Schrade", "abbas.schrade@example.com", "male", "Germany", "2025-03-01 15:01:58"),
    @(10, "Tugce", "Riezebos", "tugce.riezebos@example.com", "female", "Netherlands", "2025-03-01 15:01:58")
)

$row = 7
foreach ($record in $data) {
    $ws.Cells.Item($row, 1).Value = $record[0]
    $ws.Cells.Item($row, 2).Value = $record[1]
    $ws.Cells.Item($row, 3).Value = $record[2]
    $ws.Cells.Item($row, 4).Value = $record[3]
    $ws.Cells.Item($row, 5).Value = $record[4]
    $ws.Cells.Item($row, 6).Value = $record[5]
    $ws.Cells.Item($row, 7).Value = $record[6]
    $row++
}
